$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Durations_psplib")

$ws.Range("C1").Value = 1586.869649648666

$ws.Range("A4").Value = 304042.7785
$ws.Range("B4").Value = 306298
$ws.Range("F4").Value = 102244.466
$ws.Range("G4").Value = 103074
